# Update task tracking va task trong tet
# Adds a new "23/01/04" attendance row (row 11) on the "Attendent check" sheet,
# fills in the red/yellow status cells for the "21/01/14" row (row 10), and adds
# a new "23/01/14" deadline/usecase row (row 10) on the "Deadline" sheet. Also
# moves the active sheet/selection to mirror the authored workbook state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Attendent check")
$ws2 = $wb.Worksheets.Item("Deadline")

# ---------------------------------------------------------------------------
# Sheet "Attendent check": fill in the status cells for the existing
# "21/01/14" row (row 10), and add a brand new "23/01/04" row (row 11).
# ---------------------------------------------------------------------------

# Row 10 (date 21/01/14 already present in A10) gets the same kind of
# status-cell formatting used on row 9 (B9=red/s2, C9=red/s2, D9=yellow/s4,
# E9=yellow/s4) -- for row 10 the pattern is red, theme-fill, red, red.
$ws1.Range("B9").Copy() | Out-Null
$ws1.Range("B10").PasteSpecial(-4122) | Out-Null
$ws1.Range("D6").Copy() | Out-Null
$ws1.Range("C10").PasteSpecial(-4122) | Out-Null
$ws1.Range("B9").Copy() | Out-Null
$ws1.Range("D10").PasteSpecial(-4122) | Out-Null
$ws1.Range("B9").Copy() | Out-Null
$ws1.Range("E10").PasteSpecial(-4122) | Out-Null

# Row 11: brand new tracking date "23/01/04" plus its status cells.
$ws1.Range("A11").Value = "23/01/04"
$ws1.Range("B6").Copy() | Out-Null
$ws1.Range("B11").PasteSpecial(-4122) | Out-Null
$ws1.Range("B6").Copy() | Out-Null
$ws1.Range("C11").PasteSpecial(-4122) | Out-Null
$ws1.Range("C6").Copy() | Out-Null
$ws1.Range("D11").PasteSpecial(-4122) | Out-Null
$ws1.Range("D6").Copy() | Out-Null
$ws1.Range("E11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "Deadline": add a new "23/01/14" deadline row (row 10), matching the
# formatting of the existing "19/01/14 08:00" row (row 9).
# ---------------------------------------------------------------------------

$ws2.Range("A9").Copy() | Out-Null
$ws2.Range("A10").PasteSpecial(-4122) | Out-Null
$ws2.Range("B9").Copy() | Out-Null
$ws2.Range("B10").PasteSpecial(-4122) | Out-Null
$ws2.Range("C9").Copy() | Out-Null
$ws2.Range("C10").PasteSpecial(-4122) | Out-Null
$ws2.Range("D9").Copy() | Out-Null
$ws2.Range("D10").PasteSpecial(-4122) | Out-Null
$ws2.Range("E9").Copy() | Out-Null
$ws2.Range("E10").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

$ws2.Range("A10").Value = "23/01/14"
$ws2.Range("B10").Value = "Usecase"
$ws2.Range("C10").Value = "Usecase"
$ws2.Range("D10").Value = "Usecase"
$ws2.Range("E10").Value = "Usecase"

# ---------------------------------------------------------------------------
# Selection / active sheet bookkeeping, matching the authored workbook state:
# "Attendent check" ends up with E11 selected (no longer the active tab) while
# "Deadline" becomes the active tab with E10 selected.
# ---------------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("E11").Select()

$ws2.Activate()
$ws2.Range("E10").Select()
